# The "H 72" data row (worksheet row 2) was removed from the table.
# Removing it shifts every row below it up by one, so the last row
# (formerly row 63, "SC 232") disappears and the sheet's used range
# shrinks from A1:F63 to A1:F62.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Delete()
